$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 575734
$ws.Range("R2").Value = 6671439

$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
